# Apply updated simulation results to Sheet1, range B2:M25
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = New-Object 'object[,]' 24,12
$data[0,0] = 0.7028063186220379
$data[0,1] = 0.04515230317632302
$data[0,2] = 0.1696043453175058
$data[0,3] = 0.065064173581451
$data[0,4] = 2.882730028648638
$data[0,5] = 0
$data[0,6] = 0.07973214163530429
$data[0,7] = 2.379696603339113
$data[0,8] = 0
$data[0,9] = 0.5539231946162886
$data[0,10] = 0.2350104351937574
$data[0,11] = 0.2051997460671195
$data[1,0] = 0.6830582863522352
$data[1,1] = 0.03931581342193624
$data[1,2] = 0.1687664912116489
$data[1,3] = 0.06508989977272428
$data[1,4] = 2.829799947878499
$data[1,5] = 0
$data[1,6] = 0.07973214163530429
$data[1,7] = 2.346872366106723
$data[1,8] = 0
$data[1,9] = 0.5277790340457784
$data[1,10] = 0.2320207262118359
$data[1,11] = 0.2009107545716695
$data[2,0] = 0.6714924377060072
$data[2,1] = 0.03572706087450683
$data[2,2] = 0.1682345253062287
$data[2,3] = 0.06511933386729929
$data[2,4] = 2.798130741212816
$data[2,5] = 0
$data[2,6] = 0.07973214163530429
$data[2,7] = 2.327252651222707
$data[2,8] = 0
$data[2,9] = 0.512176839428264
$data[2,10] = 0.2302950901070488
$data[2,11] = 0.1984151525754179
$data[3,0] = 0.6669200062992786
$data[3,1] = 0.03426322171019081
$data[3,2] = 0.1680133123278793
$data[3,3] = 0.06513476492838244
$data[3,4] = 2.785433527912559
$data[3,5] = 0
$data[3,6] = 0.07973214163530429
$data[3,7] = 2.31939152825089
$data[3,8] = 0
$data[3,9] = 0.5059319503062767
$data[3,10] = 0.2296195664454288
$data[3,11] = 0.1974328564106287
$data[4,0] = 0.6661692598586342
$data[4,1] = 0.03402006512612843
$data[4,2] = 0.1679763115719268
$data[4,3] = 0.06513753500359698
$data[4,4] = 2.783337728820044
$data[4,5] = 0
$data[4,6] = 0.07973214163530429
$data[4,7] = 2.318094282575359
$data[4,8] = 0
$data[4,9] = 0.5049018204340285
$data[4,10] = 0.2295090688687296
$data[4,11] = 0.1972718421135795
$data[5,0] = 0.671430202353946
$data[5,1] = 0.03570732480692129
$data[5,2] = 0.1682315599340178
$data[5,3] = 0.06511952805251298
$data[5,4] = 2.797958659629927
$data[5,5] = 0
$data[5,6] = 0.07973214163530429
$data[5,7] = 2.327146090978374
$data[5,8] = 0
$data[5,9] = 0.5120921607721129
$data[5,10] = 0.2302858676375124
$data[5,11] = 0.1984017645367757
$data[6,0] = 0.6958810868925411
$data[6,1] = 0.04314088021632756
$data[6,2] = 0.1693190732270082
$data[6,3] = 0.06507021750912045
$data[6,4] = 2.864307142575512
$data[6,5] = 0
$data[6,6] = 0.07973214163530429
$data[6,7] = 2.368267657414179
$data[6,8] = 0
$data[6,9] = 0.5448151113477877
$data[6,10] = 0.2339567468037203
$data[6,11] = 0.2036922974812825
$data[7,0] = 0.7482707620662552
$data[7,1] = 0.05768279880936689
$data[7,2] = 0.1713139159676835
$data[7,3] = 0.06508145503366158
$data[7,4] = 3.001032914592372
$data[7,5] = 0
$data[7,6] = 0.07973214163530429
$data[7,7] = 2.453168246310227
$data[7,8] = 0
$data[7,9] = 0.6125691691672159
$data[7,10] = 0.2420287271484938
$data[7,11] = 0.2151610416290204
$data[8,0] = 0.7894780231037544
$data[8,1] = 0.06835368576710721
$data[8,2] = 0.1726973073640607
$data[8,3] = 0.06515518555748656
$data[8,4] = 3.105573187011913
$data[8,5] = 0
$data[8,6] = 0.07973214163530429
$data[8,7] = 2.518178197096063
$data[8,8] = 0
$data[8,9] = 0.6645544679179238
$data[8,10] = 0.2484928594395228
$data[8,11] = 0.2242557369543619
$data[9,0] = 0.8088165435355563
$data[9,1] = 0.07320719269830533
$data[9,2] = 0.1733091725649913
$data[9,3] = 0.06520287962822024
$data[9,4] = 3.15403135467659
$data[9,5] = 0
$data[9,6] = 0.07973214163530429
$data[9,7] = 2.548332896438595
$data[9,8] = 0
$data[9,9] = 0.6886880951509795
$data[9,10] = 0.251549771719823
$data[9,11] = 0.2285387979917388
$data[10,0] = 0.8162249136177877
$data[10,1] = 0.07504511545373305
$data[10,2] = 0.1735383896688987
$data[10,3] = 0.06522296934319094
$data[10,4] = 3.172511697295846
$data[10,5] = 0
$data[10,6] = 0.07973214163530429
$data[10,7] = 2.559835791264419
$data[10,8] = 0
$data[10,9] = 0.6978969355419622
$data[10,10] = 0.2527240837254254
$data[10,11] = 0.2301816658636611
$data[11,0] = 0.8146255950983061
$data[11,1] = 0.07464928391399894
$data[11,2] = 0.1734891336860755
$data[11,3] = 0.06521855252243469
$data[11,4] = 3.168525821689059
$data[11,5] = 0
$data[11,6] = 0.07973214163530429
$data[11,7] = 2.557354696823282
$data[11,8] = 0
$data[11,9] = 0.6959105325816779
$data[11,10] = 0.2524704308130339
$data[11,11] = 0.2298269124990284
$data[12,0] = 0.8094243251928219
$data[12,1] = 0.07335839897326935
$data[12,2] = 0.1733280800094334
$data[12,3] = 0.06520449179262044
$data[12,4] = 3.155549129347406
$data[12,5] = 0
$data[12,6] = 0.07973214163530429
$data[12,7] = 2.54927756163093
$data[12,8] = 0
$data[12,9] = 0.6894443088412459
$data[12,10] = 0.2516460478924358
$data[12,11] = 0.2286735375664293
$data[13,0] = 0.8062495048718006
$data[13,1] = 0.07256769861817247
$data[13,2] = 0.1732291073882664
$data[13,3] = 0.06519614324091094
$data[13,4] = 3.14761751380999
$data[13,5] = 0
$data[13,6] = 0.07973214163530429
$data[13,7] = 2.54434102871717
$data[13,8] = 0
$data[13,9] = 0.6854926763596154
$data[13,10] = 0.2511432675983372
$data[13,11] = 0.2279697923104536
$data[14,0] = 0.7882261405229372
$data[14,1] = 0.06803649234572617
$data[14,2] = 0.1726569716769255
$data[14,3] = 0.06515235290806309
$data[14,4] = 3.102424533939285
$data[14,5] = 0
$data[14,6] = 0.07973214163530429
$data[14,7] = 2.516219246147898
$data[14,8] = 0
$data[14,9] = 0.6629870584596631
$data[14,10] = 0.2482954237893154
$data[14,11] = 0.2239787629603285
$data[15,0] = 0.7773213036313393
$data[15,1] = 0.06525663462832654
$data[15,2] = 0.1723015349190469
$data[15,3] = 0.06512910956831774
$data[15,4] = 3.074931581649338
$data[15,5] = 0
$data[15,6] = 0.07973214163530429
$data[15,7] = 2.499116601288549
$data[15,8] = 0
$data[15,9] = 0.6493050024437821
$data[15,10] = 0.2465781592439242
$data[15,11] = 0.2215677439115353
$data[16,0] = 0.7711049577071378
$data[16,1] = 0.06365765834253523
$data[16,2] = 0.1720954534200239
$data[16,3] = 0.06511707319752524
$data[16,4] = 3.059203256364697
$data[16,5] = 0
$data[16,6] = 0.07973214163530429
$data[16,7] = 2.489334309342482
$data[16,8] = 0
$data[16,9] = 0.6414811089138936
$data[16,10] = 0.2456013853792172
$data[16,11] = 0.2201947167987797
$data[17,0] = 0.7690097983575299
$data[17,1] = 0.0631162561019778
$data[17,2] = 0.1720253946717456
$data[17,3] = 0.06511322695408239
$data[17,4] = 3.053892479840783
$data[17,5] = 0
$data[17,6] = 0.07973214163530429
$data[17,7] = 2.486031576457577
$data[17,8] = 0
$data[17,9] = 0.6388399103021527
$data[17,10] = 0.2452725475822177
$data[17,11] = 0.2197321912809187
$data[18,0] = 0.7784763641891743
$data[18,1] = 0.06555256225821893
$data[18,2] = 0.172339541694992
$data[18,3] = 0.06513144598971721
$data[18,4] = 3.077849463305711
$data[18,5] = 0
$data[18,6] = 0.07973214163530429
$data[18,7] = 2.500931543460197
$data[18,8] = 0
$data[18,9] = 0.6507567534102634
$data[18,10] = 0.2467598316074344
$data[18,11] = 0.2218229803319574
$data[19,0] = 0.810949749404358
$data[19,1] = 0.07373756215199023
$data[19,2] = 0.1733754525605136
$data[19,3] = 0.06520856675289011
$data[19,4] = 3.159357158189607
$data[19,5] = 0
$data[19,6] = 0.07973214163530429
$data[19,7] = 2.551647730313874
$data[19,8] = 0
$data[19,9] = 0.6913416953441924
$data[19,10] = 0.2518877352437272
$data[19,11] = 0.2290117427810401
$data[20,0] = 0.8326701364687779
$data[20,1] = 0.07908701380296179
$data[20,2] = 0.1740380193969244
$data[20,3] = 0.06527079335998387
$data[20,4] = 3.213386810096011
$data[20,5] = 0
$data[20,6] = 0.07973214163530429
$data[20,7] = 2.585283272280051
$data[20,8] = 0
$data[20,9] = 0.718274167234398
$data[20,10] = 0.2553366050599521
$data[20,11] = 0.2338322165362428
$data[21,0] = 0.821032065696528
$data[21,1] = 0.07623186666229742
$data[21,2] = 0.1736857096933946
$data[21,3] = 0.0652365020332546
$data[21,4] = 3.184480479338106
$data[21,5] = 0
$data[21,6] = 0.07973214163530429
$data[21,7] = 2.567286431778086
$data[21,8] = 0
$data[21,9] = 0.7038624189367226
$data[21,10] = 0.2534869595044711
$data[21,11] = 0.231248259041152
$data[22,0] = 0.7779539962556044
$data[22,1] = 0.06541877587470424
$data[22,2] = 0.1723223642415377
$data[22,3] = 0.06513038556141382
$data[22,4] = 3.076530046894078
$data[22,5] = 0
$data[22,6] = 0.07973214163530429
$data[22,7] = 2.500110851613996
$data[22,8] = 0
$data[22,9] = 0.6501002856699643
$data[22,10] = 0.2466776648128786
$data[22,11] = 0.221707547142131
$data[23,0] = 0.7336215740657792
$data[23,1] = 0.05375173029666769
$data[23,2] = 0.1707888236846209
$data[23,3] = 0.06506689370168495
$data[23,4] = 2.963331030722543
$data[23,5] = 0
$data[23,6] = 0.07973214163530429
$data[23,7] = 2.429740502187997
$data[23,8] = 0
$data[23,9] = 0.5938539922609323
$data[23,10] = 0.2397514180237295
$data[23,11] = 0.2119411528005344

$range = $ws.Range("B2:M25")
$range.Value2 = $data
